# Insert a new row before row 334, shifting existing rows 334-408 down to 335-409,
# then populate the new row 334 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 334 (pushes old rows 334..408 to 335..409)
$ws.Rows.Item(334).Insert()

# Populate the newly inserted row 334 with its data
$ws.Cells.Item(334, 1).Value2  = 5
$ws.Cells.Item(334, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(334, 3).Value2  = "Maule"
$ws.Cells.Item(334, 4).Value2  = 44798
$ws.Cells.Item(334, 5).Value2  = 7
$ws.Cells.Item(334, 6).Value2  = 100112032
$ws.Cells.Item(334, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(334, 8).Value2  = "Sin especificar"
$ws.Cells.Item(334, 9).Value2  = "Primera"
$ws.Cells.Item(334, 10).Value2 = 300
$ws.Cells.Item(334, 11).Value2 = 22000
$ws.Cells.Item(334, 12).Value2 = 22000
$ws.Cells.Item(334, 13).Value2 = 22000
$ws.Cells.Item(334, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(334, 15).Value2 = "Región del Maule"
$ws.Cells.Item(334, 16).Value2 = 440
$ws.Cells.Item(334, 17).Value2 = 50
$ws.Cells.Item(334, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the date-formatted style used by the other rows
$ws.Cells.Item(334, 4).NumberFormat = $ws.Cells.Item(335, 4).NumberFormat
